$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The first column (A) is removed entirely; every subsequent column
# shifts one position to the left (B->A, C->B, D->C, E->D, F->E).
$ws.Range("A1").EntireColumn.Delete()
